$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("aquaculture_production")
$ws2 = $wb.Worksheets.Item("wild_caught")

# --- 1. Add the new "Sheet1" worksheet at the end, matching the diff's new sheet ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws3.Name = "Sheet1"

# Populate A1:G46 with the salmon time-series snapshot (values only, no formulas)
$arr = New-Object 'object[,]' 46,7
$arr[0,0] = "year"
$arr[0,1] = "BC_FarmedAtlantic"
$arr[0,2] = "BC_FarmedChinook"
$arr[0,3] = "BC_FarmedCoho"
$arr[0,4] = "BC_TotalFarmedSalmon"
$arr[0,5] = "Canada_FarmedSalmon"
$arr[0,6] = "Canada_FarmedAtlantic"
$arr[1,0] = 1977
$arr[1,6] = 0
$arr[2,0] = 1978
$arr[2,6] = 0
$arr[3,0] = 1979
$arr[3,6] = 5
$arr[4,0] = 1980
$arr[4,6] = 27
$arr[5,0] = 1981
$arr[5,6] = 76
$arr[6,0] = 1982
$arr[6,6] = 143
$arr[7,0] = 1983
$arr[7,6] = 68
$arr[8,0] = 1984
$arr[8,6] = 222
$arr[9,0] = 1985
$arr[9,6] = 349
$arr[10,0] = 1986
$arr[10,1] = 0
$arr[10,2] = 87
$arr[10,3] = 304
$arr[10,6] = 682
$arr[11,0] = 1987
$arr[11,1] = 3
$arr[11,2] = 949
$arr[11,3] = 791
$arr[11,6] = 1385
$arr[12,0] = 1988
$arr[12,1] = 80
$arr[12,2] = 3545
$arr[12,3] = 2743
$arr[12,6] = 3431
$arr[13,0] = 1989
$arr[13,1] = 1280
$arr[13,2] = 8514
$arr[13,3] = 1815
$arr[13,6] = 5967
$arr[14,0] = 1990
$arr[14,1] = 1640
$arr[14,2] = 10396
$arr[14,3] = 1296
$arr[14,6] = 9625
$arr[15,0] = 1991
$arr[15,4] = 24362
$arr[15,5] = 34109
$arr[15,6] = 13499
$arr[16,0] = 1992
$arr[16,4] = 19814
$arr[16,5] = 30325
$arr[16,6] = 17305
$arr[17,0] = 1993
$arr[17,4] = 25555
$arr[17,5] = 36670
$arr[17,6] = 23483
$arr[18,0] = 1994
$arr[18,4] = 23657
$arr[18,5] = 23657
$arr[18,6] = 27773
$arr[19,0] = 1995
$arr[19,4] = 27275
$arr[19,5] = 42515
$arr[19,6] = 33674
$arr[20,0] = 1996
$arr[20,4] = 27756
$arr[20,5] = 45624
$arr[20,6] = 36475
$arr[21,0] = 1997
$arr[21,4] = 36465
$arr[21,5] = 56775
$arr[21,6] = 51015
$arr[22,0] = 1998
$arr[22,4] = 42200
$arr[22,5] = 58618
$arr[22,6] = 49475
$arr[23,0] = 1999
$arr[23,4] = 49700
$arr[23,5] = 72890
$arr[23,6] = 61990
$arr[24,0] = 2000
$arr[24,4] = 49000
$arr[24,5] = 82195
$arr[24,6] = 72495
$arr[25,0] = 2001
$arr[25,4] = 68000
$arr[25,5] = 105606
$arr[25,6] = 95606
$arr[26,0] = 2002
$arr[26,4] = 84200
$arr[26,5] = 126321
$arr[26,6] = 114921
$arr[27,0] = 2003
$arr[27,4] = 65411
$arr[27,5] = 99961
$arr[27,6] = 107228
$arr[28,0] = 2004
$arr[28,4] = 55646
$arr[28,5] = 90646
$arr[28,6] = 96774
$arr[29,0] = 2005
$arr[29,4] = 63370
$arr[29,5] = 98370
$arr[29,6] = 98370
$arr[30,0] = 2006
$arr[30,4] = 70181
$arr[30,5] = 118061
$arr[30,6] = 118061
$arr[31,0] = 2007
$arr[31,4] = 70998
$arr[31,5] = 102509
$arr[31,6] = 102509
$arr[32,0] = 2008
$arr[32,4] = 73265
$arr[32,5] = 104075
$arr[32,6] = 104075
$arr[33,0] = 2009
$arr[33,4] = 68662
$arr[33,5] = 100212
$arr[33,6] = 100212
$arr[34,0] = 2010
$arr[34,4] = 70831
$arr[34,5] = 101544
$arr[34,6] = 101544
$arr[35,0] = 2011
$arr[35,4] = 83144
$arr[35,5] = 110328
$arr[35,6] = 110328
$arr[36,0] = 2012
$arr[36,4] = 79981
$arr[36,5] = 116101
$arr[36,6] = 116101
$arr[37,0] = 2013
$arr[37,4] = 74673
$arr[37,5] = 100027
$arr[37,6] = 97629
$arr[38,0] = 2014
$arr[38,4] = 54971
$arr[38,5] = 78979
$arr[38,6] = 86347
$arr[39,0] = 2015
$arr[39,4] = 92926
$arr[39,5] = 121926
$arr[39,6] = 121926
$arr[40,0] = 2016
$arr[40,4] = 90511
$arr[40,5] = 123522
$arr[40,6] = 123522
$arr[41,0] = 2017
$arr[41,4] = 85608
$arr[41,5] = 120553
$arr[41,6] = 120553
$arr[42,0] = 2018
$arr[42,4] = 87010
$arr[42,5] = 123184
$arr[42,6] = 123184
$arr[43,0] = 2019
$arr[43,4] = 88874
$arr[43,5] = 118632
$arr[43,6] = 118630
$arr[44,0] = 2020
$arr[44,4] = 91666
$arr[44,5] = 120285
$arr[44,6] = 120427
$arr[45,0] = 2021
$arr[45,4] = 84171
$arr[45,5] = 120186
$arr[45,6] = 120186

$rng = $ws3.Range($ws3.Cells.Item(1,1), $ws3.Cells.Item(46,7))
$rng.Value = $arr

# Column widths for the new sheet
$ws3.Columns.Item(2).ColumnWidth = 22
$ws3.Columns.Item(3).ColumnWidth = 18.33203125
$ws3.Columns.Item(4).ColumnWidth = 18.83203125
$ws3.Columns.Item(5).ColumnWidth = 22.33203125
$ws3.Columns.Item(6).ColumnWidth = 22.33203125
$ws3.Columns.Item(7).ColumnWidth = 28.5

# F31 keeps the distinct font style that it already carried on aquaculture_production
$ws3.Range("F31").Font.Color = $ws1.Range("F31").Font.Color

# --- 2. Convert the per-row SUM formulas on aquaculture_production!F into shared formula groups ---
$ws1.Range("F16:F30").Formula = "=SUM(E16:E16)"
$ws1.Range("F32:F37").Formula = "=SUM(E32:E32)"
$ws1.Range("F39:F46").Formula = "=SUM(E39:E39)"

# --- 3. Update view/selection state: wild_caught becomes the active/selected sheet ---
$ws1.Range("F46").Select()

# Leave the new sheet's own selection/view as captured in the diff
$ws3.Activate()
$ws3.Range("D8").Select()

$ws2.Activate()
$ws2.Range("B76").Select()
